$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New file identities being reported on in this handoff-status run:
#   14501b0b-4f8d-4432-b987-389a73f11de7.png   (was 03ba9452-...-bb223b0e6fad.md)
#   b2a55f6c-c600-4713-8ae1-13d8fb59b505.png   (new)
#   d563eebd-3af1-4301-b896-2e2d0e8feea7.md    (new)
# ---------------------------------------------------------------------------

$overviewDate = "2016-50-17 18:50:08"
$zhDate       = "2016-03-17 18:50:00"
$deDate       = "2016-03-17 18:50:08"
$epoch        = "0001-01-01 00:00:00"
$ready        = "Ready for handoff"

# =============================== Overview ==================================
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Hyperlinks.Delete()

$ws1.Range("B2").Value = $ready
$ws1.Range("C2").Value = $ready
$ws1.Range("D2").Value = $overviewDate

$ws1.Range("B3").Value = $ready
$ws1.Range("C3").Value = $ready
$ws1.Range("D3").Value = $overviewDate

$ws1.Range("B4").Value = $ready
$ws1.Range("C4").Value = $ready
$ws1.Range("D4").Value = $overviewDate

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/85d73a499ebcf4d4e2b485ddca015266d4358ad7/e2e/14501b0b-4f8d-4432-b987-389a73f11de7.png", "", "", "14501b0b-4f8d-4432-b987-389a73f11de7.png")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/85d73a499ebcf4d4e2b485ddca015266d4358ad7/e2e/b2a55f6c-c600-4713-8ae1-13d8fb59b505.png", "", "", "b2a55f6c-c600-4713-8ae1-13d8fb59b505.png")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/85d73a499ebcf4d4e2b485ddca015266d4358ad7/e2e/d563eebd-3af1-4301-b896-2e2d0e8feea7.md", "", "", "d563eebd-3af1-4301-b896-2e2d0e8feea7.md")

# ================================ zh-cn =====================================
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Delete()

$ws2.Range("C2").Value = $ready
$ws2.Range("E2").Value = $zhDate
$ws2.Range("H2").Value = $epoch
$ws2.Range("I2").Value = "IsDependency"
$ws2.Range("J2").Value = "e2e\d563eebd-3af1-4301-b896-2e2d0e8feea7.md"

$ws2.Range("C3").Value = $ready
$ws2.Range("E3").Value = $zhDate
$ws2.Range("H3").Value = $epoch
$ws2.Range("I3").Value = "IsDependency"
$ws2.Range("J3").Value = "e2e\d563eebd-3af1-4301-b896-2e2d0e8feea7.md"

$ws2.Range("C4").Value = $ready
$ws2.Range("E4").Value = $zhDate
$ws2.Range("H4").Value = $epoch
$ws2.Range("I4").Value = "Include"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/85d73a499ebcf4d4e2b485ddca015266d4358ad7/e2e/14501b0b-4f8d-4432-b987-389a73f11de7.png", "", "", "14501b0b-4f8d-4432-b987-389a73f11de7.png")
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/85d73a499ebcf4d4e2b485ddca015266d4358ad7/e2e/14501b0b-4f8d-4432-b987-389a73f11de7.png", "", "", ".png")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7fbdbc6d1b7a91b538666b503d3bfcc2501a44d9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b8cab2530e923f9576a50685746fa2282f47aecc.png", "", "", "b8cab2530e923f9576a50685746fa2282f47aecc.png")

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/85d73a499ebcf4d4e2b485ddca015266d4358ad7/e2e/b2a55f6c-c600-4713-8ae1-13d8fb59b505.png", "", "", "b2a55f6c-c600-4713-8ae1-13d8fb59b505.png")
$ws2.Hyperlinks.Add($ws2.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/85d73a499ebcf4d4e2b485ddca015266d4358ad7/e2e/b2a55f6c-c600-4713-8ae1-13d8fb59b505.png", "", "", ".png")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7fbdbc6d1b7a91b538666b503d3bfcc2501a44d9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a6f100277d48058d79a84e5b4f46b0eca4f9541f.png", "", "", "a6f100277d48058d79a84e5b4f46b0eca4f9541f.png")

$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/85d73a499ebcf4d4e2b485ddca015266d4358ad7/e2e/d563eebd-3af1-4301-b896-2e2d0e8feea7.md", "", "", "d563eebd-3af1-4301-b896-2e2d0e8feea7.md")
$ws2.Hyperlinks.Add($ws2.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/85d73a499ebcf4d4e2b485ddca015266d4358ad7/e2e/d563eebd-3af1-4301-b896-2e2d0e8feea7.md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7fbdbc6d1b7a91b538666b503d3bfcc2501a44d9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d563eebd-3af1-4301-b896-2e2d0e8feea7.f5b479cdc079ef60dde4688e6ad08650807a9dbe.zh-cn.xlf", "", "", "d563eebd-3af1-4301-b896-2e2d0e8feea7.f5b479cdc079ef60dde4688e6ad08650807a9dbe.zh-cn.xlf")

# ================================ de-de =====================================
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Delete()

$ws3.Range("C2").Value = $ready
$ws3.Range("E2").Value = $deDate
$ws3.Range("H2").Value = $epoch
$ws3.Range("I2").Value = "IsDependency"
$ws3.Range("J2").Value = "e2e\d563eebd-3af1-4301-b896-2e2d0e8feea7.md"

$ws3.Range("C3").Value = $ready
$ws3.Range("E3").Value = $deDate
$ws3.Range("H3").Value = $epoch
$ws3.Range("I3").Value = "IsDependency"
$ws3.Range("J3").Value = "e2e\d563eebd-3af1-4301-b896-2e2d0e8feea7.md"

$ws3.Range("C4").Value = $ready
$ws3.Range("E4").Value = $deDate
$ws3.Range("H4").Value = $epoch
$ws3.Range("I4").Value = "Include"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/85d73a499ebcf4d4e2b485ddca015266d4358ad7/e2e/14501b0b-4f8d-4432-b987-389a73f11de7.png", "", "", "14501b0b-4f8d-4432-b987-389a73f11de7.png")
$ws3.Hyperlinks.Add($ws3.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/85d73a499ebcf4d4e2b485ddca015266d4358ad7/e2e/14501b0b-4f8d-4432-b987-389a73f11de7.png", "", "", ".png")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/054a3590de57a2adabaaca3961be477a028625fe/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b8cab2530e923f9576a50685746fa2282f47aecc.png", "", "", "b8cab2530e923f9576a50685746fa2282f47aecc.png")

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/85d73a499ebcf4d4e2b485ddca015266d4358ad7/e2e/b2a55f6c-c600-4713-8ae1-13d8fb59b505.png", "", "", "b2a55f6c-c600-4713-8ae1-13d8fb59b505.png")
$ws3.Hyperlinks.Add($ws3.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/85d73a499ebcf4d4e2b485ddca015266d4358ad7/e2e/b2a55f6c-c600-4713-8ae1-13d8fb59b505.png", "", "", ".png")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/054a3590de57a2adabaaca3961be477a028625fe/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a6f100277d48058d79a84e5b4f46b0eca4f9541f.png", "", "", "a6f100277d48058d79a84e5b4f46b0eca4f9541f.png")

$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/85d73a499ebcf4d4e2b485ddca015266d4358ad7/e2e/d563eebd-3af1-4301-b896-2e2d0e8feea7.md", "", "", "d563eebd-3af1-4301-b896-2e2d0e8feea7.md")
$ws3.Hyperlinks.Add($ws3.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/85d73a499ebcf4d4e2b485ddca015266d4358ad7/e2e/d563eebd-3af1-4301-b896-2e2d0e8feea7.md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/054a3590de57a2adabaaca3961be477a028625fe/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d563eebd-3af1-4301-b896-2e2d0e8feea7.f5b479cdc079ef60dde4688e6ad08650807a9dbe.de-de.xlf", "", "", "d563eebd-3af1-4301-b896-2e2d0e8feea7.f5b479cdc079ef60dde4688e6ad08650807a9dbe.de-de.xlf")
